$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 3000
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 3000
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 9000
$ws.Cells.Item(29, 14).Value = -9562

$ws.Cells.Item(38, 8).Value = 896.7143
$ws.Cells.Item(38, 9).Value = 455.4
$ws.Cells.Item(38, 10).Value = 2000
$ws.Cells.Item(38, 11).Value = 1366.2
$ws.Cells.Item(38, 12).Value = 6000
$ws.Cells.Item(38, 13).Value = -994.1999999999998
$ws.Cells.Item(38, 14).Value = -6744

$ws.Cells.Item(92, 8).Value = 1537.5625
$ws.Cells.Item(92, 9).Value = 1536.909
$ws.Cells.Item(92, 10).Value = 1539
$ws.Cells.Item(92, 11).Value = 1536.909
$ws.Cells.Item(92, 12).Value = 1539
$ws.Cells.Item(92, 13).Value = -288.9090000000001

$ws.Cells.Item(96, 8).Value = 1632.125
$ws.Cells.Item(96, 9).Value = 428
$ws.Cells.Item(96, 10).Value = 5244.5
$ws.Cells.Item(96, 11).Value = 1284
$ws.Cells.Item(96, 12).Value = 15733.5
$ws.Cells.Item(96, 13).Value = 89

$ws.Cells.Item(97, 8).Value = 1447
$ws.Cells.Item(97, 9).Value = 800
$ws.Cells.Item(97, 10).Value = 1511.7
$ws.Cells.Item(97, 11).Value = 2400
$ws.Cells.Item(97, 12).Value = 4535.1
$ws.Cells.Item(97, 13).Value = -1904
$ws.Cells.Item(97, 14).Value = -5527.1

$ws.Cells.Item(101, 8).Value = 985
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 985
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 2955
$ws.Cells.Item(101, 14).Value = -6199

$ws.Cells.Item(103, 8).Value = 1847.7646
$ws.Cells.Item(103, 9).Value = 1519
$ws.Cells.Item(103, 10).Value = 2027.091
$ws.Cells.Item(103, 11).Value = 4557
$ws.Cells.Item(103, 12).Value = 6081.272999999999
$ws.Cells.Item(103, 13).Value = -3971

$ws.Cells.Item(106, 8).Value = 3333
$ws.Cells.Item(106, 9).Value = 3333
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 3333
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = -2702

$ws.Cells.Item(112, 8).Value = 3364.2856
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 3364.2856
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 10092.8568
$ws.Cells.Item(112, 13).Value = ""
$ws.Cells.Item(112, 14).Value = -12308.8568

$ws.Cells.Item(125, 8).Value = 515.4
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 515.4
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 4638.599999999999
$ws.Cells.Item(125, 13).Value = ""
$ws.Cells.Item(125, 14).Value = -9558.599999999999

$ws.Cells.Item(132, 8).Value = 1657.95
$ws.Cells.Item(132, 9).Value = 1657.95
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4973.85
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2443.85

$ws.Cells.Item(135, 8).Value = 859.9091
$ws.Cells.Item(135, 9).Value = 686
$ws.Cells.Item(135, 10).Value = 2599
$ws.Cells.Item(135, 11).Value = 6174
$ws.Cells.Item(135, 12).Value = 23391
$ws.Cells.Item(135, 13).Value = -3639
$ws.Cells.Item(135, 14).Value = -28461

$ws.Cells.Item(138, 8).Value = 5692.121
$ws.Cells.Item(138, 9).Value = 1922.875
$ws.Cells.Item(138, 10).Value = 6898.28
$ws.Cells.Item(138, 11).Value = 5768.625
$ws.Cells.Item(138, 12).Value = 20694.84
$ws.Cells.Item(138, 13).Value = -628.625
$ws.Cells.Item(138, 14).Value = -30974.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12957.577
$ws.Cells.Item(32, 9).Value = 11386.826
$ws.Cells.Item(32, 10).Value = 25000
$ws.Cells.Item(32, 11).Value = 11386.826
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = -11099.826

$ws.Cells.Item(45, 8).Value = 2675.889
$ws.Cells.Item(45, 9).Value = 2675.889
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 2675.889
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -2298.889

$ws.Cells.Item(61, 8).Value = 3690.2856
$ws.Cells.Item(61, 9).Value = 3690.2856
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 3690.2856
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -3478.2856

$ws.Cells.Item(63, 8).Value = 7525
$ws.Cells.Item(63, 9).Value = 1286
$ws.Cells.Item(63, 10).Value = 20003
$ws.Cells.Item(63, 11).Value = 1286
$ws.Cells.Item(63, 12).Value = 20003
$ws.Cells.Item(63, 13).Value = -600
$ws.Cells.Item(63, 14).Value = -21375

$ws.Cells.Item(66, 8).Value = 7525
$ws.Cells.Item(66, 9).Value = 1286
$ws.Cells.Item(66, 10).Value = 20003
$ws.Cells.Item(66, 11).Value = 6430
$ws.Cells.Item(66, 12).Value = 100015
$ws.Cells.Item(66, 13).Value = -2998
$ws.Cells.Item(66, 14).Value = -106879

$ws.Cells.Item(74, 8).Value = 13666.5
$ws.Cells.Item(74, 9).Value = 12497.125
$ws.Cells.Item(74, 10).Value = 16005.25
$ws.Cells.Item(74, 11).Value = 12497.125
$ws.Cells.Item(74, 12).Value = 16005.25
$ws.Cells.Item(74, 13).Value = -11623.125

$ws.Cells.Item(77, 8).Value = 13666.5
$ws.Cells.Item(77, 9).Value = 12497.125
$ws.Cells.Item(77, 10).Value = 16005.25
$ws.Cells.Item(77, 11).Value = 62485.625
$ws.Cells.Item(77, 12).Value = 80026.25
$ws.Cells.Item(77, 13).Value = -58117.625

$ws.Cells.Item(97, 8).Value = 4459.857
$ws.Cells.Item(97, 9).Value = 649.6667
$ws.Cells.Item(97, 10).Value = 7317.5
$ws.Cells.Item(97, 11).Value = 649.6667
$ws.Cells.Item(97, 12).Value = 7317.5
$ws.Cells.Item(97, 13).Value = -153.6667

$ws.Cells.Item(132, 8).Value = 3491.3076
$ws.Cells.Item(132, 9).Value = 2861.5
$ws.Cells.Item(132, 10).Value = 4499
$ws.Cells.Item(132, 11).Value = 8584.5
$ws.Cells.Item(132, 12).Value = 13497
$ws.Cells.Item(132, 13).Value = -6054.5

$ws.Cells.Item(136, 8).Value = 3690.2856
$ws.Cells.Item(136, 9).Value = 3690.2856
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 11070.8568
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -8520.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 9874.125
$ws.Cells.Item(20, 9).Value = 15348.5
$ws.Cells.Item(20, 10).Value = 8049.3335
$ws.Cells.Item(20, 11).Value = 15348.5
$ws.Cells.Item(20, 12).Value = 8049.3335
$ws.Cells.Item(20, 13).Value = -15101.5

$ws.Cells.Item(86, 8).Value = 9040
$ws.Cells.Item(86, 9).Value = 2600
$ws.Cells.Item(86, 10).Value = 13333.333
$ws.Cells.Item(86, 11).Value = 2600
$ws.Cells.Item(86, 12).Value = 13333.333
$ws.Cells.Item(86, 13).Value = -1477

$ws.Cells.Item(89, 8).Value = 9040
$ws.Cells.Item(89, 9).Value = 2600
$ws.Cells.Item(89, 10).Value = 13333.333
$ws.Cells.Item(89, 11).Value = 13000
$ws.Cells.Item(89, 12).Value = 66666.66500000001
$ws.Cells.Item(89, 13).Value = -7384

$ws.Cells.Item(134, 8).Value = 1949.5
$ws.Cells.Item(134, 9).Value = 1949.5
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 5848.5
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -3313.5

$ws.Cells.Item(135, 8).Value = 80159.664
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 80159.664
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 80159.664
$ws.Cells.Item(135, 14).Value = -90299.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2217.3076
$ws.Cells.Item(31, 9).Value = 1982.5
$ws.Cells.Item(31, 10).Value = 3000
$ws.Cells.Item(31, 11).Value = 1982.5
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = -1687.5
$ws.Cells.Item(31, 14).Value = -3590

$ws.Cells.Item(33, 8).Value = 6809.4
$ws.Cells.Item(33, 9).Value = 5503
$ws.Cells.Item(33, 10).Value = 12035
$ws.Cells.Item(33, 11).Value = 5503
$ws.Cells.Item(33, 12).Value = 12035
$ws.Cells.Item(33, 13).Value = -5124
$ws.Cells.Item(33, 14).Value = -12793

$ws.Cells.Item(34, 8).Value = 2217.3076
$ws.Cells.Item(34, 9).Value = 1982.5
$ws.Cells.Item(34, 10).Value = 3000
$ws.Cells.Item(34, 11).Value = 1982.5
$ws.Cells.Item(34, 12).Value = 3000
$ws.Cells.Item(34, 13).Value = -1780.5
$ws.Cells.Item(34, 14).Value = -3404

$ws.Cells.Item(58, 8).Value = 4894.857
$ws.Cells.Item(58, 9).Value = 3191
$ws.Cells.Item(58, 10).Value = 7166.6665
$ws.Cells.Item(58, 11).Value = 3191
$ws.Cells.Item(58, 12).Value = 7166.6665
$ws.Cells.Item(58, 13).Value = -2988
$ws.Cells.Item(58, 14).Value = -7572.6665

$ws.Cells.Item(132, 8).Value = 3005
$ws.Cells.Item(132, 9).Value = 1995.4
$ws.Cells.Item(132, 10).Value = 4014.6
$ws.Cells.Item(132, 11).Value = 5986.200000000001
$ws.Cells.Item(132, 12).Value = 12043.8
$ws.Cells.Item(132, 13).Value = -3456.200000000001

$ws.Cells.Item(136, 8).Value = 4894.857
$ws.Cells.Item(136, 9).Value = 3191
$ws.Cells.Item(136, 10).Value = 7166.6665
$ws.Cells.Item(136, 11).Value = 9573
$ws.Cells.Item(136, 12).Value = 21499.9995
$ws.Cells.Item(136, 13).Value = -7023
$ws.Cells.Item(136, 14).Value = -26599.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 6192.857
$ws.Cells.Item(3, 9).Value = 6558.3335
$ws.Cells.Item(3, 10).Value = 4000
$ws.Cells.Item(3, 11).Value = 19675.0005
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = -19563.0005

$ws.Cells.Item(5, 8).Value = 1039.5
$ws.Cells.Item(5, 9).Value = 1594
$ws.Cells.Item(5, 10).Value = 485
$ws.Cells.Item(5, 11).Value = 4782
$ws.Cells.Item(5, 12).Value = 1455
$ws.Cells.Item(5, 13).Value = -4670
$ws.Cells.Item(5, 14).Value = -1679

$ws.Cells.Item(50, 8).Value = 100
$ws.Cells.Item(50, 9).Value = 100
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 300
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 181

$ws.Cells.Item(53, 8).Value = 100
$ws.Cells.Item(53, 9).Value = 100
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 300
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 181

$ws.Cells.Item(135, 8).Value = 1039.5
$ws.Cells.Item(135, 9).Value = 1594
$ws.Cells.Item(135, 10).Value = 485
$ws.Cells.Item(135, 11).Value = 14346
$ws.Cells.Item(135, 12).Value = 4365
$ws.Cells.Item(135, 13).Value = -11811
$ws.Cells.Item(135, 14).Value = -9435

$ws.Cells.Item(136, 8).Value = 2264.6667
$ws.Cells.Item(136, 9).Value = 2029.3334
$ws.Cells.Item(136, 10).Value = 2500
$ws.Cells.Item(136, 11).Value = 6088.0002
$ws.Cells.Item(136, 12).Value = 7500
$ws.Cells.Item(136, 13).Value = -988.0002000000004
$ws.Cells.Item(136, 14).Value = -17700

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 7499.5
$ws.Cells.Item(102, 9).Value = 7499.5
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 7499.5
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -5877.5

$ws.Cells.Item(122, 8).Value = 4180.722
$ws.Cells.Item(122, 9).Value = 2940.9375
$ws.Cells.Item(122, 10).Value = 14099
$ws.Cells.Item(122, 11).Value = 8822.8125
$ws.Cells.Item(122, 12).Value = 42297
$ws.Cells.Item(122, 13).Value = -6372.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = ""
$ws.Cells.Item(29, 14).Value = ""

$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = ""

$ws.Cells.Item(93, 8).Value = 2666.4167
$ws.Cells.Item(93, 9).Value = 2666.5
$ws.Cells.Item(93, 10).Value = 2666.3333
$ws.Cells.Item(93, 11).Value = 2666.5
$ws.Cells.Item(93, 12).Value = 2666.3333
$ws.Cells.Item(93, 13).Value = -1418.5
$ws.Cells.Item(93, 14).Value = -5162.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 25995.8
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 25995.8
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 25995.8
$ws.Cells.Item(101, 14).Value = -32485.8

$ws.Cells.Item(132, 8).Value = 4807.3076
$ws.Cells.Item(132, 9).Value = 2613.8572
$ws.Cells.Item(132, 10).Value = 7366.3335
$ws.Cells.Item(132, 11).Value = 7841.571599999999
$ws.Cells.Item(132, 12).Value = 22099.0005
$ws.Cells.Item(132, 13).Value = -5311.571599999999

$ws.Cells.Item(136, 8).Value = 1040.4286
$ws.Cells.Item(136, 9).Value = 1040.4286
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 3121.2858
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -571.2857999999997
